# Fix a typo in the "commerceVariableName" column (D) of the Tracker sheet:
# "fireDomain.insall" -> "fireDomain.install" for rows 4 and 5.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")

$ws.Range("D4").Value = "fireDomain.install"
$ws.Range("D5").Value = "fireDomain.install"

# Update the active selection to match the authored state.
$ws.Range("D8").Select()
